$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at the top of the sheet, shifting everything down by one row.
$ws.Rows.Item(1).Insert()

# Add the note about part substitution in the newly created A1 cell.
$ws.Range("A1").Value = "NOTE: Altitude on the SparkFun forum points out that the pads for the 0805 parts below are small for hand soldering and that substituting 0603 parts where possible actually makes assembly easier.  If this works for you, there is no functional downside to using the 0603 parts. Thx. - CS"

# Update the workbook-level defined name so it still covers the BOM table,
# which has shifted down by one row because of the inserted row.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!VoltageReferenceProgrammableRevDBOM") {
        $n.RefersTo = "=Sheet1!`$A`$2:`$M`$26"
    }
}
